$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.254929423332214
$ws.Range("B1").Value = 2.40125036239624
$ws.Range("C1").Value = 4.118087291717529
$ws.Range("D1").Value = 2.747961759567261
$ws.Range("E1").Value = 1.355506300926208
